# "Generate Report for Handback" - update the localization-status report:
#  - Overview sheet: status changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" (for both zh-cn and de-de columns)
#  - zh-cn / de-de sheets: Status column gets the same new text, the
#    "Latest Target File" / "Latest Handback File" columns get filled in
#    (with a hyperlink on the target-file cell), and the
#    "Latest Handback DateTime" is stamped.
#  - Column widths are widened for the columns that now hold long file
#    names / datetimes.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

$ovRepoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f3be9c27d5c1f74e381b2de93a8f8b3e28bf86a/e2e/9ea1fc09-9903-4ea7-b9a1-fd74419e2798.md"
$targetFileName = "9ea1fc09-9903-4ea7-b9a1-fd74419e2798.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $handedBack
$wsOverview.Range("F2").Value = $handedBack
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

$wsOverview.Range("E1").ColumnWidth = 29.166666666666664
$wsOverview.Range("F1").ColumnWidth = 29.166666666666664

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $handedBack
$wsZh.Range("C3").Value = $handedBack

$wsZh.Range("I2").Value = $targetFileName
$wsZh.Range("I2").Style = "Hyperlink"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $ovRepoUrl, "", "", $targetFileName)

$wsZh.Range("J2").Value = "9ea1fc09-9903-4ea7-b9a1-fd74419e2798.2c988c245430e50ab5c50bf0b0204c19e79ef1f0.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-31 01:07:42"

$wsZh.Range("I3").Value = $targetFileName
$wsZh.Range("I3").Style = "Hyperlink"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $ovRepoUrl, "", "", $targetFileName)

$wsZh.Range("J3").Value = "9ea1fc09-9903-4ea7-b9a1-fd74419e2798.2c988c245430e50ab5c50bf0b0204c19e79ef1f0.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-31 01:07:42"

$wsZh.Range("C1").ColumnWidth = 29.166666666666664
$wsZh.Range("I1").ColumnWidth = 39.16666666666667
$wsZh.Range("J1").ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $handedBack
$wsDe.Range("C3").Value = $handedBack

$wsDe.Range("I2").Value = $targetFileName
$wsDe.Range("I2").Style = "Hyperlink"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $ovRepoUrl, "", "", $targetFileName)

$wsDe.Range("J2").Value = "9ea1fc09-9903-4ea7-b9a1-fd74419e2798.2c988c245430e50ab5c50bf0b0204c19e79ef1f0.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-31 01:07:49"

$wsDe.Range("I3").Value = $targetFileName
$wsDe.Range("I3").Style = "Hyperlink"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $ovRepoUrl, "", "", $targetFileName)

$wsDe.Range("J3").Value = "9ea1fc09-9903-4ea7-b9a1-fd74419e2798.2c988c245430e50ab5c50bf0b0204c19e79ef1f0.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-31 01:07:49"

$wsDe.Range("C1").ColumnWidth = 29.166666666666664
$wsDe.Range("I1").ColumnWidth = 39.16666666666667
$wsDe.Range("J1").ColumnWidth = 39.16666666666667
